# Update "想去人数" (interested-count) figures in the 展览 (sheet1) and
# 全部类型 (sheet4) worksheets, plus a single value on 演出 (sheet2),
# to match the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value  = 243
$ws1.Cells.Item(4, 6).Value  = 238
$ws1.Cells.Item(5, 6).Value  = 1776
$ws1.Cells.Item(6, 6).Value  = 664
$ws1.Cells.Item(7, 6).Value  = 320
$ws1.Cells.Item(8, 6).Value  = 494
$ws1.Cells.Item(9, 6).Value  = 4546
$ws1.Cells.Item(10, 6).Value = 54
$ws1.Cells.Item(13, 6).Value = 983
$ws1.Cells.Item(14, 6).Value = 1284
$ws1.Cells.Item(17, 6).Value = 2977
$ws1.Cells.Item(18, 6).Value = 1807
$ws1.Cells.Item(20, 6).Value = 43
$ws1.Cells.Item(22, 6).Value = 23
$ws1.Cells.Item(24, 6).Value = 928
$ws1.Cells.Item(26, 6).Value = 31
$ws1.Cells.Item(27, 6).Value = 2343
$ws1.Cells.Item(28, 6).Value = 998
$ws1.Cells.Item(29, 6).Value = 2405
$ws1.Cells.Item(30, 6).Value = 246
$ws1.Cells.Item(31, 6).Value = 1113
$ws1.Cells.Item(32, 6).Value = 579
$ws1.Cells.Item(34, 6).Value = 883
$ws1.Cells.Item(35, 6).Value = 419
$ws1.Cells.Item(36, 6).Value = 1108
$ws1.Cells.Item(37, 6).Value = 906
$ws1.Cells.Item(38, 6).Value = 1176
$ws1.Cells.Item(40, 6).Value = 846
$ws1.Cells.Item(41, 6).Value = 521
$ws1.Cells.Item(42, 6).Value = 359
$ws1.Cells.Item(43, 6).Value = 280
$ws1.Cells.Item(44, 6).Value = 3479

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(10, 6).Value = 879

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value  = 243
$ws4.Cells.Item(4, 6).Value  = 238
$ws4.Cells.Item(6, 6).Value  = 1776
$ws4.Cells.Item(7, 6).Value  = 664
$ws4.Cells.Item(8, 6).Value  = 320
$ws4.Cells.Item(9, 6).Value  = 494
$ws4.Cells.Item(10, 6).Value = 4546
$ws4.Cells.Item(11, 6).Value = 54
$ws4.Cells.Item(15, 6).Value = 1284
$ws4.Cells.Item(16, 6).Value = 2977
$ws4.Cells.Item(18, 6).Value = 1807
$ws4.Cells.Item(20, 6).Value = 43
$ws4.Cells.Item(23, 6).Value = 879
$ws4.Cells.Item(25, 6).Value = 23
$ws4.Cells.Item(26, 6).Value = 928
$ws4.Cells.Item(28, 6).Value = 2343
$ws4.Cells.Item(31, 6).Value = 998
$ws4.Cells.Item(33, 6).Value = 2405
$ws4.Cells.Item(34, 6).Value = 1113
$ws4.Cells.Item(35, 6).Value = 579
$ws4.Cells.Item(36, 6).Value = 883
$ws4.Cells.Item(37, 6).Value = 1108
$ws4.Cells.Item(38, 6).Value = 906
$ws4.Cells.Item(40, 6).Value = 1176
$ws4.Cells.Item(41, 6).Value = 846
$ws4.Cells.Item(42, 6).Value = 521
$ws4.Cells.Item(44, 6).Value = 359
$ws4.Cells.Item(47, 6).Value = 280
$ws4.Cells.Item(48, 6).Value = 3479
